$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price record is inserted right above the existing
# row 70, pushing every following row down by one (old row 70 -> 71, ...,
# old row 169 -> 170). Insert a whole new row so the rest of the table
# shifts automatically.
$ws.Rows.Item(70).EntireRow.Insert()

# Populate the newly inserted row 70 with the new "Choclo - Choclero"
# record for Comercializadora del Agro de Limarí.
$ws.Cells.Item(70, 1).Value = 2
$ws.Cells.Item(70, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44994
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112024
$ws.Cells.Item(70, 7).Value = "Choclo"
$ws.Cells.Item(70, 8).Value = "Choclero"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 30000
$ws.Cells.Item(70, 11).Value = 350
$ws.Cells.Item(70, 12).Value = 400
$ws.Cells.Item(70, 13).Value = 375
$ws.Cells.Item(70, 14).Value = "$/unidad"
$ws.Cells.Item(70, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(70, 16).Value = 375
$ws.Cells.Item(70, 17).Value = 1
$ws.Cells.Item(70, 18).Value = "Hortaliza"
